$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (Matteo Mazzola's team) with the new teammates
$ws.Range("C2").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("D2").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("E2").Value = "marco bertolini | Fc Wanda Tim"
$ws.Range("F2").Value = "Matteo  Tatarella | Bayern Mona"

# Remove rows 3 and 4 entirely (old "Raffa prova" and "Nicola Lorenzi" teams)
$ws.Range("A3:F4").EntireRow.Delete()
